$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 387.5  # H18: 400 -> 387.5
$ws.Cells.Item(18, 10).Value = 350  # J18: 0 -> 350
$ws.Cells.Item(18, 12).Value = 350  # L18: 0 -> 350
$ws.Cells.Item(18, 14).Value = -918  # N18: None -> -918
$ws.Cells.Item(46, 8).Value = 4895  # H46: 4897 -> 4895
$ws.Cells.Item(46, 9).Value = 0  # I46: 4898 -> 0
$ws.Cells.Item(46, 10).Value = 4895  # J46: 4895.5 -> 4895
$ws.Cells.Item(46, 11).Value = 0  # K46: 14694 -> 0
$ws.Cells.Item(46, 12).Value = 14685  # L46: 14686.5 -> 14685
$ws.Cells.Item(46, 13).ClearContents()  # M46: was -14575
$ws.Cells.Item(46, 14).Value = -14923  # N46: -14924.5 -> -14923
$ws.Cells.Item(60, 8).Value = 4895  # H60: 4897 -> 4895
$ws.Cells.Item(60, 9).Value = 0  # I60: 4898 -> 0
$ws.Cells.Item(60, 10).Value = 4895  # J60: 4895.5 -> 4895
$ws.Cells.Item(60, 11).Value = 0  # K60: 14694 -> 0
$ws.Cells.Item(60, 12).Value = 14685  # L60: 14686.5 -> 14685
$ws.Cells.Item(60, 13).ClearContents()  # M60: was -14210
$ws.Cells.Item(60, 14).Value = -15653  # N60: -15654.5 -> -15653
$ws.Cells.Item(74, 8).Value = 4279  # H74: 4382.6665 -> 4279
$ws.Cells.Item(74, 9).Value = 4100  # I74: 4060.2 -> 4100
$ws.Cells.Item(74, 10).Value = 4995  # J74: 5995 -> 4995
$ws.Cells.Item(74, 11).Value = 4100  # K74: 4060.2 -> 4100
$ws.Cells.Item(74, 12).Value = 4995  # L74: 5995 -> 4995
$ws.Cells.Item(74, 13).Value = -3164  # M74: -3124.2 -> -3164
$ws.Cells.Item(74, 14).Value = -6867  # N74: -7867 -> -6867
$ws.Cells.Item(76, 8).Value = 6521.125  # H76: 6881.4287 -> 6521.125
$ws.Cells.Item(76, 9).Value = 6717.2  # I76: 7396.75 -> 6717.2
$ws.Cells.Item(76, 11).Value = 6717.2  # K76: 7396.75 -> 6717.2
$ws.Cells.Item(76, 13).Value = -6402.2  # M76: -7081.75 -> -6402.2
$ws.Cells.Item(77, 8).Value = 4279  # H77: 4382.6665 -> 4279
$ws.Cells.Item(77, 9).Value = 4100  # I77: 4060.2 -> 4100
$ws.Cells.Item(77, 10).Value = 4995  # J77: 5995 -> 4995
$ws.Cells.Item(77, 11).Value = 20500  # K77: 20301 -> 20500
$ws.Cells.Item(77, 12).Value = 24975  # L77: 29975 -> 24975
$ws.Cells.Item(77, 13).Value = -15820  # M77: -15621 -> -15820
$ws.Cells.Item(77, 14).Value = -34335  # N77: -39335 -> -34335
$ws.Cells.Item(79, 8).Value = 6521.125  # H79: 6881.4287 -> 6521.125
$ws.Cells.Item(79, 9).Value = 6717.2  # I79: 7396.75 -> 6717.2
$ws.Cells.Item(79, 11).Value = 6717.2  # K79: 7396.75 -> 6717.2
$ws.Cells.Item(79, 13).Value = -5625.2  # M79: -6304.75 -> -5625.2
$ws.Cells.Item(80, 8).Value = 38722.938  # H80: 38734.688 -> 38722.938
$ws.Cells.Item(80, 9).Value = 75922.25  # I80: 75945.75 -> 75922.25
$ws.Cells.Item(80, 11).Value = 227766.75  # K80: 227837.25 -> 227766.75
$ws.Cells.Item(80, 13).Value = -226768.75  # M80: -226839.25 -> -226768.75
$ws.Cells.Item(83, 8).Value = 38722.938  # H83: 38734.688 -> 38722.938
$ws.Cells.Item(83, 9).Value = 75922.25  # I83: 75945.75 -> 75922.25
$ws.Cells.Item(83, 11).Value = 683300.25  # K83: 683511.75 -> 683300.25
$ws.Cells.Item(83, 13).Value = -678308.25  # M83: -678519.75 -> -678308.25
$ws.Cells.Item(100, 8).Value = 2605.4688  # H100: 2311.0466 -> 2605.4688
$ws.Cells.Item(100, 9).Value = 2099.2693  # I100: 1907.5946 -> 2099.2693
$ws.Cells.Item(100, 11).Value = 2099.2693  # K100: 1907.5946 -> 2099.2693
$ws.Cells.Item(100, 13).Value = -1558.2693  # M100: -1366.5946 -> -1558.2693
$ws.Cells.Item(111, 8).Value = 1849  # H111: 1908.8 -> 1849
$ws.Cells.Item(111, 9).Value = 1857.3334  # I111: 1936.25 -> 1857.3334
$ws.Cells.Item(111, 11).Value = 5572.0002  # K111: 5808.75 -> 5572.0002
$ws.Cells.Item(111, 13).Value = -2505.0002  # M111: -2741.75 -> -2505.0002
$ws.Cells.Item(113, 8).Value = 4061.8572  # H113: 4806.4546 -> 4061.8572
$ws.Cells.Item(113, 9).Value = 2941.111  # I113: 3745.8333 -> 2941.111
$ws.Cells.Item(113, 11).Value = 2941.111  # K113: 3745.8333 -> 2941.111
$ws.Cells.Item(113, 13).Value = 312.8890000000001  # M113: -491.8332999999998 -> 312.8890000000001
$ws.Cells.Item(116, 8).Value = 275545.2  # H116: 3243.3333 -> 275545.2
$ws.Cells.Item(116, 9).Value = 502266.66  # I116: 2931.6667 -> 502266.66
$ws.Cells.Item(116, 10).Value = 3479.4  # J116: 3399.1667 -> 3479.4
$ws.Cells.Item(116, 11).Value = 502266.66  # K116: 2931.6667 -> 502266.66
$ws.Cells.Item(116, 12).Value = 3479.4  # L116: 3399.1667 -> 3479.4
$ws.Cells.Item(116, 13).Value = -498824.66  # M116: 510.3332999999998 -> -498824.66
$ws.Cells.Item(116, 14).Value = -10363.4  # N116: -10283.1667 -> -10363.4
$ws.Cells.Item(125, 8).Value = 1552.5  # H125: 1754.5714 -> 1552.5
$ws.Cells.Item(125, 9).Value = 760.3333  # I125: 974.5 -> 760.3333
$ws.Cells.Item(125, 10).Value = 2344.6667  # J125: 2066.6 -> 2344.6667
$ws.Cells.Item(125, 11).Value = 6842.9997  # K125: 8770.5 -> 6842.9997
$ws.Cells.Item(125, 12).Value = 21102.0003  # L125: 18599.4 -> 21102.0003
$ws.Cells.Item(125, 13).Value = -4382.9997  # M125: -6310.5 -> -4382.9997
$ws.Cells.Item(125, 14).Value = -26022.0003  # N125: -23519.4 -> -26022.0003
$ws.Cells.Item(127, 8).Value = 2998.25  # H127: 3049.5 -> 2998.25
$ws.Cells.Item(127, 9).Value = 3999.5  # I127: 3399.6667 -> 3999.5
$ws.Cells.Item(127, 10).Value = 1997  # J127: 1999 -> 1997
$ws.Cells.Item(127, 11).Value = 11998.5  # K127: 10199.0001 -> 11998.5
$ws.Cells.Item(127, 12).Value = 5991  # L127: 5997 -> 5991
$ws.Cells.Item(127, 13).Value = -7038.5  # M127: -5239.000100000001 -> -7038.5
$ws.Cells.Item(127, 14).Value = -15911  # N127: -15917 -> -15911
$ws.Cells.Item(131, 8).Value = 3069.5  # H131: 1952.8125 -> 3069.5
$ws.Cells.Item(131, 9).Value = 604.0909  # I131: 615 -> 604.0909
$ws.Cells.Item(131, 10).Value = 6082.778  # J131: 7750 -> 6082.778
$ws.Cells.Item(131, 11).Value = 1812.2727  # K131: 1845 -> 1812.2727
$ws.Cells.Item(131, 12).Value = 18248.334  # L131: 23250 -> 18248.334
$ws.Cells.Item(131, 13).Value = 3227.7273  # M131: 3195 -> 3227.7273
$ws.Cells.Item(131, 14).Value = -28328.334  # N131: -33330 -> -28328.334
$ws.Cells.Item(132, 8).Value = 6230.76  # H132: 8047.222 -> 6230.76
$ws.Cells.Item(132, 9).Value = 5001.579  # I132: 5082.1577 -> 5001.579
$ws.Cells.Item(132, 10).Value = 10123.167  # J132: 15089.25 -> 10123.167
$ws.Cells.Item(132, 11).Value = 15004.737  # K132: 15246.4731 -> 15004.737
$ws.Cells.Item(132, 12).Value = 30369.501  # L132: 45267.75 -> 30369.501
$ws.Cells.Item(132, 13).Value = -12474.737  # M132: -12716.4731 -> -12474.737
$ws.Cells.Item(132, 14).Value = -35429.501  # N132: -50327.75 -> -35429.501
$ws.Cells.Item(138, 8).Value = 4715.224  # H138: 4643.6 -> 4715.224
$ws.Cells.Item(138, 10).Value = 3726.239  # J138: 3677.9167 -> 3726.239
$ws.Cells.Item(138, 12).Value = 11178.717  # L138: 11033.7501 -> 11178.717
$ws.Cells.Item(138, 14).Value = -21458.717  # N138: -21313.7501 -> -21458.717
$ws.Cells.Item(141, 8).Value = 13093.134  # H141: 13743 -> 13093.134
$ws.Cells.Item(141, 9).Value = 13093.134  # I141: 13743 -> 13093.134
$ws.Cells.Item(141, 11).Value = 39279.402  # K141: 41229 -> 39279.402
$ws.Cells.Item(141, 13).Value = -34099.402  # M141: -36049 -> -34099.402

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 882.9394  # H2: 973.1667 -> 882.9394
$ws.Cells.Item(2, 9).Value = 879.3461  # I2: 952.7083 -> 879.3461
$ws.Cells.Item(2, 10).Value = 896.2857  # J2: 1055 -> 896.2857
$ws.Cells.Item(2, 11).Value = 879.3461  # K2: 952.7083 -> 879.3461
$ws.Cells.Item(2, 12).Value = 896.2857  # L2: 1055 -> 896.2857
$ws.Cells.Item(2, 13).Value = -766.3461  # M2: -839.7083 -> -766.3461
$ws.Cells.Item(2, 14).Value = -1122.2857  # N2: -1281 -> -1122.2857
$ws.Cells.Item(32, 8).Value = 6325.9155  # H32: 9252.584000000001 -> 6325.9155
$ws.Cells.Item(32, 9).Value = 4913.8506  # I32: 5120.6562 -> 4913.8506
$ws.Cells.Item(32, 10).Value = 29978  # J32: 29594.385 -> 29978
$ws.Cells.Item(32, 11).Value = 4913.8506  # K32: 5120.6562 -> 4913.8506
$ws.Cells.Item(32, 12).Value = 29978  # L32: 29594.385 -> 29978
$ws.Cells.Item(32, 13).Value = -4626.8506  # M32: -4833.6562 -> -4626.8506
$ws.Cells.Item(32, 14).Value = -30552  # N32: -30168.385 -> -30552
$ws.Cells.Item(45, 8).Value = 912606.9399999999  # H45: 837465.4399999999 -> 912606.9399999999
$ws.Cells.Item(45, 9).Value = 1668496.6  # I45: 1668831.6 -> 1668496.6
$ws.Cells.Item(45, 10).Value = 5539.2  # J45: 6099.1665 -> 5539.2
$ws.Cells.Item(45, 11).Value = 1668496.6  # K45: 1668831.6 -> 1668496.6
$ws.Cells.Item(45, 12).Value = 5539.2  # L45: 6099.1665 -> 5539.2
$ws.Cells.Item(45, 13).Value = -1668119.6  # M45: -1668454.6 -> -1668119.6
$ws.Cells.Item(45, 14).Value = -6293.2  # N45: -6853.1665 -> -6293.2
$ws.Cells.Item(110, 8).Value = 1500.1  # H110: 1646.7142 -> 1500.1
$ws.Cells.Item(110, 9).Value = 1361.2222  # I110: 1530.579 -> 1361.2222
$ws.Cells.Item(110, 11).Value = 1361.2222  # K110: 1530.579 -> 1361.2222
$ws.Cells.Item(110, 13).Value = 683.7778000000001  # M110: 514.421 -> 683.7778000000001
$ws.Cells.Item(116, 8).Value = 882.9394  # H116: 973.1667 -> 882.9394
$ws.Cells.Item(116, 9).Value = 879.3461  # I116: 952.7083 -> 879.3461
$ws.Cells.Item(116, 10).Value = 896.2857  # J116: 1055 -> 896.2857
$ws.Cells.Item(116, 11).Value = 879.3461  # K116: 952.7083 -> 879.3461
$ws.Cells.Item(116, 12).Value = 896.2857  # L116: 1055 -> 896.2857
$ws.Cells.Item(116, 13).Value = 1414.6539  # M116: 1341.2917 -> 1414.6539
$ws.Cells.Item(116, 14).Value = -5484.2857  # N116: -5643 -> -5484.2857
$ws.Cells.Item(132, 8).Value = 47638.25  # H132: 58907.05 -> 47638.25
$ws.Cells.Item(132, 9).Value = 60186.055  # I132: 81481.92 -> 60186.055
$ws.Cells.Item(132, 11).Value = 180558.165  # K132: 244445.76 -> 180558.165
$ws.Cells.Item(132, 13).Value = -178028.165  # M132: -241915.76 -> -178028.165

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 882.9394  # H3: 973.1667 -> 882.9394
$ws.Cells.Item(3, 9).Value = 879.3461  # I3: 952.7083 -> 879.3461
$ws.Cells.Item(3, 10).Value = 896.2857  # J3: 1055 -> 896.2857
$ws.Cells.Item(3, 11).Value = 879.3461  # K3: 952.7083 -> 879.3461
$ws.Cells.Item(3, 12).Value = 896.2857  # L3: 1055 -> 896.2857
$ws.Cells.Item(3, 13).Value = -765.3461  # M3: -838.7083 -> -765.3461
$ws.Cells.Item(3, 14).Value = -1124.2857  # N3: -1283 -> -1124.2857
$ws.Cells.Item(20, 8).Value = 2056.2104  # H20: 1984.9565 -> 2056.2104
$ws.Cells.Item(20, 9).Value = 2007.1333  # I20: 1929.1177 -> 2007.1333
$ws.Cells.Item(20, 10).Value = 2240.25  # J20: 2143.1667 -> 2240.25
$ws.Cells.Item(20, 11).Value = 2007.1333  # K20: 1929.1177 -> 2007.1333
$ws.Cells.Item(20, 12).Value = 2240.25  # L20: 2143.1667 -> 2240.25
$ws.Cells.Item(20, 13).Value = -1760.1333  # M20: -1682.1177 -> -1760.1333
$ws.Cells.Item(20, 14).Value = -2734.25  # N20: -2637.1667 -> -2734.25
$ws.Cells.Item(54, 8).Value = 5499.3335  # H54: 5666.3335 -> 5499.3335
$ws.Cells.Item(54, 9).Value = 5499.3335  # I54: 5666.3335 -> 5499.3335
$ws.Cells.Item(54, 11).Value = 5499.3335  # K54: 5666.3335 -> 5499.3335
$ws.Cells.Item(54, 13).Value = -5015.3335  # M54: -5182.3335 -> -5015.3335
$ws.Cells.Item(86, 8).Value = 3633.0557  # H86: 3647.4736 -> 3633.0557
$ws.Cells.Item(86, 9).Value = 3366.5  # I86: 3408 -> 3366.5
$ws.Cells.Item(86, 10).Value = 4166.1665  # J86: 4166.3335 -> 4166.1665
$ws.Cells.Item(86, 11).Value = 3366.5  # K86: 3408 -> 3366.5
$ws.Cells.Item(86, 12).Value = 4166.1665  # L86: 4166.3335 -> 4166.1665
$ws.Cells.Item(86, 13).Value = -2243.5  # M86: -2285 -> -2243.5
$ws.Cells.Item(86, 14).Value = -6412.1665  # N86: -6412.3335 -> -6412.1665
$ws.Cells.Item(89, 8).Value = 3633.0557  # H89: 3647.4736 -> 3633.0557
$ws.Cells.Item(89, 9).Value = 3366.5  # I89: 3408 -> 3366.5
$ws.Cells.Item(89, 10).Value = 4166.1665  # J89: 4166.3335 -> 4166.1665
$ws.Cells.Item(89, 11).Value = 16832.5  # K89: 17040 -> 16832.5
$ws.Cells.Item(89, 12).Value = 20830.8325  # L89: 20831.6675 -> 20830.8325
$ws.Cells.Item(89, 13).Value = -11216.5  # M89: -11424 -> -11216.5
$ws.Cells.Item(89, 14).Value = -32062.8325  # N89: -32063.6675 -> -32062.8325
$ws.Cells.Item(94, 8).Value = 1074.3889  # H94: 1042.6 -> 1074.3889
$ws.Cells.Item(94, 9).Value = 994.8889  # I94: 875.6667 -> 994.8889
$ws.Cells.Item(94, 11).Value = 994.8889  # K94: 875.6667 -> 994.8889
$ws.Cells.Item(94, 13).Value = -543.8889  # M94: -424.6667 -> -543.8889
$ws.Cells.Item(99, 8).Value = 1619  # H99: 1704.5555 -> 1619
$ws.Cells.Item(99, 9).Value = 1521.1111  # I99: 1605.125 -> 1521.1111
$ws.Cells.Item(99, 11).Value = 1521.1111  # K99: 1605.125 -> 1521.1111
$ws.Cells.Item(99, 13).Value = -23.11110000000008  # M99: -107.125 -> -23.11110000000008
$ws.Cells.Item(105, 8).Value = 2048.0715  # H105: 1919.7222 -> 2048.0715
$ws.Cells.Item(105, 9).Value = 1697.1666  # I105: 1710.8572 -> 1697.1666
$ws.Cells.Item(105, 10).Value = 2311.25  # J105: 2052.6365 -> 2311.25
$ws.Cells.Item(105, 11).Value = 1697.1666  # K105: 1710.8572 -> 1697.1666
$ws.Cells.Item(105, 12).Value = 2311.25  # L105: 2052.6365 -> 2311.25
$ws.Cells.Item(105, 13).Value = 49.83339999999998  # M105: 36.14280000000008 -> 49.83339999999998
$ws.Cells.Item(105, 14).Value = -5805.25  # N105: -5546.636500000001 -> -5805.25
$ws.Cells.Item(115, 8).Value = 84999.336  # H115: 85000 -> 84999.336
$ws.Cells.Item(115, 10).Value = 84999.336  # J115: 85000 -> 84999.336
$ws.Cells.Item(115, 12).Value = 84999.336  # L115: 85000 -> 84999.336
$ws.Cells.Item(115, 14).Value = -88133.336  # N115: -88134 -> -88133.336
$ws.Cells.Item(132, 8).Value = 99719.75  # H132: 110397.4 -> 99719.75
$ws.Cells.Item(132, 10).Value = 99719.75  # J132: 110397.4 -> 99719.75
$ws.Cells.Item(132, 12).Value = 99719.75  # L132: 110397.4 -> 99719.75
$ws.Cells.Item(132, 14).Value = -109839.75  # N132: -120517.4 -> -109839.75
$ws.Cells.Item(134, 8).Value = 3453.4707  # H134: 3492.2285 -> 3453.4707
$ws.Cells.Item(134, 9).Value = 3363.0938  # I134: 3363.4375 -> 3363.0938
$ws.Cells.Item(134, 10).Value = 4899.5  # J134: 4866 -> 4899.5
$ws.Cells.Item(134, 11).Value = 10089.2814  # K134: 10090.3125 -> 10089.2814
$ws.Cells.Item(134, 12).Value = 14698.5  # L134: 14598 -> 14698.5
$ws.Cells.Item(134, 13).Value = -7554.2814  # M134: -7555.3125 -> -7554.2814
$ws.Cells.Item(134, 14).Value = -19768.5  # N134: -19668 -> -19768.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 870.3913  # H22: 955.2381 -> 870.3913
$ws.Cells.Item(22, 9).Value = 955.6429000000001  # I22: 1020 -> 955.6429000000001
$ws.Cells.Item(22, 10).Value = 737.7778  # J22: 850 -> 737.7778
$ws.Cells.Item(22, 11).Value = 955.6429000000001  # K22: 1020 -> 955.6429000000001
$ws.Cells.Item(22, 12).Value = 737.7778  # L22: 850 -> 737.7778
$ws.Cells.Item(22, 13).Value = -605.6429000000001  # M22: -670 -> -605.6429000000001
$ws.Cells.Item(22, 14).Value = -1437.7778  # N22: -1550 -> -1437.7778
$ws.Cells.Item(31, 8).Value = 2775.55  # H31: 2547.9412 -> 2775.55
$ws.Cells.Item(31, 9).Value = 2395.3157  # I31: 2547.9412 -> 2395.3157
$ws.Cells.Item(31, 10).Value = 10000  # J31: 0 -> 10000
$ws.Cells.Item(31, 11).Value = 2395.3157  # K31: 2547.9412 -> 2395.3157
$ws.Cells.Item(31, 12).Value = 10000  # L31: 0 -> 10000
$ws.Cells.Item(31, 13).Value = -2100.3157  # M31: -2252.9412 -> -2100.3157
$ws.Cells.Item(31, 14).Value = -10590  # N31: None -> -10590
$ws.Cells.Item(34, 8).Value = 2775.55  # H34: 2547.9412 -> 2775.55
$ws.Cells.Item(34, 9).Value = 2395.3157  # I34: 2547.9412 -> 2395.3157
$ws.Cells.Item(34, 10).Value = 10000  # J34: 0 -> 10000
$ws.Cells.Item(34, 11).Value = 2395.3157  # K34: 2547.9412 -> 2395.3157
$ws.Cells.Item(34, 12).Value = 10000  # L34: 0 -> 10000
$ws.Cells.Item(34, 13).Value = -2193.3157  # M34: -2345.9412 -> -2193.3157
$ws.Cells.Item(34, 14).Value = -10404  # N34: None -> -10404
$ws.Cells.Item(62, 8).Value = 18998.691  # H62: 23499.9 -> 18998.691
$ws.Cells.Item(62, 10).Value = 25776  # J62: 36666.668 -> 25776
$ws.Cells.Item(62, 12).Value = 25776  # L62: 36666.668 -> 25776
$ws.Cells.Item(62, 14).Value = -27024  # N62: -37914.668 -> -27024
$ws.Cells.Item(65, 8).Value = 18998.691  # H65: 23499.9 -> 18998.691
$ws.Cells.Item(65, 10).Value = 25776  # J65: 36666.668 -> 25776
$ws.Cells.Item(65, 12).Value = 128880  # L65: 183333.34 -> 128880
$ws.Cells.Item(65, 14).Value = -135120  # N65: -189573.34 -> -135120
$ws.Cells.Item(122, 8).Value = 1685.25  # H122: 1699 -> 1685.25
$ws.Cells.Item(122, 9).Value = 1640.4286  # I122: 1766 -> 1640.4286
$ws.Cells.Item(122, 10).Value = 1999  # J122: 1498 -> 1999
$ws.Cells.Item(122, 11).Value = 4921.2858  # K122: 5298 -> 4921.2858
$ws.Cells.Item(122, 12).Value = 5997  # L122: 4494 -> 5997
$ws.Cells.Item(122, 13).Value = -2471.2858  # M122: -2848 -> -2471.2858
$ws.Cells.Item(122, 14).Value = -10897  # N122: -9394 -> -10897
$ws.Cells.Item(132, 8).Value = 2435.6365  # H132: 2312.1765 -> 2435.6365
$ws.Cells.Item(132, 9).Value = 2612.3914  # I132: 2356.4075 -> 2612.3914
$ws.Cells.Item(132, 10).Value = 2029.1  # J132: 2141.5715 -> 2029.1
$ws.Cells.Item(132, 11).Value = 7837.174199999999  # K132: 7069.2225 -> 7837.174199999999
$ws.Cells.Item(132, 12).Value = 6087.299999999999  # L132: 6424.7145 -> 6087.299999999999
$ws.Cells.Item(132, 13).Value = -5307.174199999999  # M132: -4539.2225 -> -5307.174199999999
$ws.Cells.Item(132, 14).Value = -11147.3  # N132: -11484.7145 -> -11147.3
$ws.Cells.Item(141, 8).Value = 414946.6  # H141: 403767 -> 414946.6
$ws.Cells.Item(141, 10).Value = 441585.66  # J141: 427884.8 -> 441585.66
$ws.Cells.Item(141, 12).Value = 441585.66  # L141: 427884.8 -> 441585.66
$ws.Cells.Item(141, 14).Value = -451945.66  # N141: -438244.8 -> -451945.66

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 3243.889  # H5: 3243.9443 -> 3243.889
$ws.Cells.Item(5, 9).Value = 1591.3334  # I5: 1591.6666 -> 1591.3334
$ws.Cells.Item(5, 11).Value = 4774.0002  # K5: 4774.9998 -> 4774.0002
$ws.Cells.Item(5, 13).Value = -4662.0002  # M5: -4662.9998 -> -4662.0002
$ws.Cells.Item(25, 8).Value = 1832.6666  # H25: 1874.5 -> 1832.6666
$ws.Cells.Item(25, 9).Value = 1750  # I25: 1833 -> 1750
$ws.Cells.Item(25, 10).Value = 1998  # J25: 1999 -> 1998
$ws.Cells.Item(25, 11).Value = 5250  # K25: 5499 -> 5250
$ws.Cells.Item(25, 12).Value = 5994  # L25: 5997 -> 5994
$ws.Cells.Item(25, 13).Value = -5081  # M25: -5330 -> -5081
$ws.Cells.Item(25, 14).Value = -6332  # N25: -6335 -> -6332
$ws.Cells.Item(30, 8).Value = 1832.6666  # H30: 1874.5 -> 1832.6666
$ws.Cells.Item(30, 9).Value = 1750  # I30: 1833 -> 1750
$ws.Cells.Item(30, 10).Value = 1998  # J30: 1999 -> 1998
$ws.Cells.Item(30, 11).Value = 5250  # K30: 5499 -> 5250
$ws.Cells.Item(30, 12).Value = 5994  # L30: 5997 -> 5994
$ws.Cells.Item(30, 13).Value = -5148  # M30: -5397 -> -5148
$ws.Cells.Item(30, 14).Value = -6198  # N30: -6201 -> -6198
$ws.Cells.Item(46, 8).Value = 999.5  # H46: 998.75 -> 999.5
$ws.Cells.Item(46, 9).Value = 0  # I46: 998 -> 0
$ws.Cells.Item(46, 10).Value = 999.5  # J46: 999 -> 999.5
$ws.Cells.Item(46, 11).Value = 0  # K46: 2994 -> 0
$ws.Cells.Item(46, 12).Value = 2998.5  # L46: 2997 -> 2998.5
$ws.Cells.Item(46, 13).ClearContents()  # M46: was -2903
$ws.Cells.Item(46, 14).Value = -3180.5  # N46: -3179 -> -3180.5
$ws.Cells.Item(56, 8).Value = 6508.3335  # H56: 6620.5557 -> 6508.3335
$ws.Cells.Item(56, 9).Value = 6508.3335  # I56: 6620.5557 -> 6508.3335
$ws.Cells.Item(56, 11).Value = 6508.3335  # K56: 6620.5557 -> 6508.3335
$ws.Cells.Item(56, 13).Value = -5978.3335  # M56: -6090.5557 -> -5978.3335
$ws.Cells.Item(68, 8).Value = 931.4  # H68: 942.6667 -> 931.4
$ws.Cells.Item(68, 9).Value = 999  # I68: 956.5 -> 999
$ws.Cells.Item(68, 10).Value = 830  # J68: 915 -> 830
$ws.Cells.Item(68, 11).Value = 2997  # K68: 2869.5 -> 2997
$ws.Cells.Item(68, 12).Value = 2490  # L68: 2745 -> 2490
$ws.Cells.Item(68, 13).Value = -2186  # M68: -2058.5 -> -2186
$ws.Cells.Item(68, 14).Value = -4112  # N68: -4367 -> -4112
$ws.Cells.Item(71, 8).Value = 931.4  # H71: 942.6667 -> 931.4
$ws.Cells.Item(71, 9).Value = 999  # I71: 956.5 -> 999
$ws.Cells.Item(71, 10).Value = 830  # J71: 915 -> 830
$ws.Cells.Item(71, 11).Value = 8991  # K71: 8608.5 -> 8991
$ws.Cells.Item(71, 12).Value = 7470  # L71: 8235 -> 7470
$ws.Cells.Item(71, 13).Value = -4935  # M71: -4552.5 -> -4935
$ws.Cells.Item(71, 14).Value = -15582  # N71: -16347 -> -15582
$ws.Cells.Item(118, 8).Value = 5789.125  # H118: 6093 -> 5789.125
$ws.Cells.Item(118, 9).Value = 5789.125  # I118: 6093 -> 5789.125
$ws.Cells.Item(118, 11).Value = 17367.375  # K118: 18279 -> 17367.375
$ws.Cells.Item(118, 13).Value = -16124.375  # M118: -17036 -> -16124.375
$ws.Cells.Item(119, 8).Value = 3434  # H119: 3813.0833 -> 3434
$ws.Cells.Item(119, 9).Value = 2159.7693  # I119: 2341.6365 -> 2159.7693
$ws.Cells.Item(119, 11).Value = 6479.3079  # K119: 7024.9095 -> 6479.3079
$ws.Cells.Item(119, 13).Value = -1641.3079  # M119: -2186.9095 -> -1641.3079
$ws.Cells.Item(120, 8).Value = 18330.666  # H120: 13656.833 -> 18330.666
$ws.Cells.Item(120, 9).Value = 18330.666  # I120: 13656.833 -> 18330.666
$ws.Cells.Item(120, 11).Value = 54991.99800000001  # K120: 40970.499 -> 54991.99800000001
$ws.Cells.Item(120, 13).Value = -50153.99800000001  # M120: -36132.499 -> -50153.99800000001
$ws.Cells.Item(121, 8).Value = 3027.4  # H121: 2689.5 -> 3027.4
$ws.Cells.Item(121, 9).Value = 0  # I121: 1000 -> 0
$ws.Cells.Item(121, 11).Value = 0  # K121: 3000 -> 0
$ws.Cells.Item(121, 13).ClearContents()  # M121: was -1690
$ws.Cells.Item(135, 8).Value = 3243.889  # H135: 3243.9443 -> 3243.889
$ws.Cells.Item(135, 9).Value = 1591.3334  # I135: 1591.6666 -> 1591.3334
$ws.Cells.Item(135, 11).Value = 14322.0006  # K135: 14324.9994 -> 14322.0006
$ws.Cells.Item(135, 13).Value = -11787.0006  # M135: -11789.9994 -> -11787.0006

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(14, 8).Value = 19875.75  # H14: 22499.75 -> 19875.75
$ws.Cells.Item(14, 9).Value = 9752  # I14: 15000 -> 9752
$ws.Cells.Item(14, 11).Value = 9752  # K14: 15000 -> 9752
$ws.Cells.Item(14, 13).Value = -9584  # M14: -14832 -> -9584
$ws.Cells.Item(86, 8).Value = 42283  # H86: 40052 -> 42283
$ws.Cells.Item(86, 10).Value = 42283  # J86: 40052 -> 42283
$ws.Cells.Item(86, 12).Value = 42283  # L86: 40052 -> 42283
$ws.Cells.Item(86, 14).Value = -44655  # N86: -42424 -> -44655
$ws.Cells.Item(89, 8).Value = 42283  # H89: 40052 -> 42283
$ws.Cells.Item(89, 10).Value = 42283  # J89: 40052 -> 42283
$ws.Cells.Item(89, 12).Value = 126849  # L89: 120156 -> 126849
$ws.Cells.Item(89, 14).Value = -138705  # N89: -132012 -> -138705
$ws.Cells.Item(97, 8).Value = 1361  # H97: 1385.1333 -> 1361
$ws.Cells.Item(97, 10).Value = 1225.4286  # J97: 1263.1666 -> 1225.4286
$ws.Cells.Item(97, 12).Value = 1225.4286  # L97: 1263.1666 -> 1225.4286
$ws.Cells.Item(97, 14).Value = -2217.4286  # N97: -2255.1666 -> -2217.4286
$ws.Cells.Item(102, 8).Value = 2373.6  # H102: 2221 -> 2373.6
$ws.Cells.Item(102, 9).Value = 1967  # I102: 2221 -> 1967
$ws.Cells.Item(102, 10).Value = 4000  # J102: 0 -> 4000
$ws.Cells.Item(102, 11).Value = 1967  # K102: 2221 -> 1967
$ws.Cells.Item(102, 12).Value = 4000  # L102: 0 -> 4000
$ws.Cells.Item(102, 13).Value = -345  # M102: -599 -> -345
$ws.Cells.Item(102, 14).Value = -7244  # N102: None -> -7244
$ws.Cells.Item(107, 8).Value = 56559.055  # H107: 67591.53 -> 56559.055
$ws.Cells.Item(107, 9).Value = 91438.17999999999  # I107: 100572.5 -> 91438.17999999999
$ws.Cells.Item(107, 10).Value = 1749  # J107: 1629.6 -> 1749
$ws.Cells.Item(107, 11).Value = 91438.17999999999  # K107: 100572.5 -> 91438.17999999999
$ws.Cells.Item(107, 12).Value = 1749  # L107: 1629.6 -> 1749
$ws.Cells.Item(107, 13).Value = -89518.17999999999  # M107: -98652.5 -> -89518.17999999999
$ws.Cells.Item(107, 14).Value = -5589  # N107: -5469.6 -> -5589
$ws.Cells.Item(122, 8).Value = 4442.6895  # H122: 7513.5186 -> 4442.6895
$ws.Cells.Item(122, 9).Value = 2894.5454  # I122: 3233.375 -> 2894.5454
$ws.Cells.Item(122, 10).Value = 5388.778  # J122: 9315.684999999999 -> 5388.778
$ws.Cells.Item(122, 11).Value = 8683.636200000001  # K122: 9700.125 -> 8683.636200000001
$ws.Cells.Item(122, 12).Value = 16166.334  # L122: 27947.055 -> 16166.334
$ws.Cells.Item(122, 13).Value = -6233.636200000001  # M122: -7250.125 -> -6233.636200000001
$ws.Cells.Item(122, 14).Value = -21066.334  # N122: -32847.055 -> -21066.334
$ws.Cells.Item(132, 8).Value = 91840.63  # H132: 125988.875 -> 91840.63
$ws.Cells.Item(132, 9).Value = 143560  # I132: 200699.4 -> 143560
$ws.Cells.Item(132, 10).Value = 1331.75  # J132: 1471.3334 -> 1331.75
$ws.Cells.Item(132, 11).Value = 430680  # K132: 602098.2 -> 430680
$ws.Cells.Item(132, 12).Value = 3995.25  # L132: 4414.0002 -> 3995.25
$ws.Cells.Item(132, 13).Value = -428150  # M132: -599568.2 -> -428150
$ws.Cells.Item(132, 14).Value = -9055.25  # N132: -9474.0002 -> -9055.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 736.13794  # H16: 707.71875 -> 736.13794
$ws.Cells.Item(16, 9).Value = 513.0454999999999  # I16: 503.36 -> 513.0454999999999
$ws.Cells.Item(16, 10).Value = 1437.2858  # J16: 1437.5714 -> 1437.2858
$ws.Cells.Item(16, 11).Value = 513.0454999999999  # K16: 503.36 -> 513.0454999999999
$ws.Cells.Item(16, 12).Value = 1437.2858  # L16: 1437.5714 -> 1437.2858
$ws.Cells.Item(16, 13).Value = -343.0454999999999  # M16: -333.36 -> -343.0454999999999
$ws.Cells.Item(16, 14).Value = -1777.2858  # N16: -1777.5714 -> -1777.2858
$ws.Cells.Item(22, 8).Value = 3291.4  # H22: 3362 -> 3291.4
$ws.Cells.Item(22, 9).Value = 1445.4  # I22: 1704.1666 -> 1445.4
$ws.Cells.Item(22, 10).Value = 4214.4  # J22: 4127.154 -> 4214.4
$ws.Cells.Item(22, 11).Value = 1445.4  # K22: 1704.1666 -> 1445.4
$ws.Cells.Item(22, 12).Value = 4214.4  # L22: 4127.154 -> 4214.4
$ws.Cells.Item(22, 13).Value = -1150.4  # M22: -1409.1666 -> -1150.4
$ws.Cells.Item(22, 14).Value = -4804.4  # N22: -4717.154 -> -4804.4
$ws.Cells.Item(27, 8).Value = 3291.4  # H27: 3362 -> 3291.4
$ws.Cells.Item(27, 9).Value = 1445.4  # I27: 1704.1666 -> 1445.4
$ws.Cells.Item(27, 10).Value = 4214.4  # J27: 4127.154 -> 4214.4
$ws.Cells.Item(27, 11).Value = 1445.4  # K27: 1704.1666 -> 1445.4
$ws.Cells.Item(27, 12).Value = 4214.4  # L27: 4127.154 -> 4214.4
$ws.Cells.Item(27, 13).Value = -1338.4  # M27: -1597.1666 -> -1338.4
$ws.Cells.Item(27, 14).Value = -4428.4  # N27: -4341.154 -> -4428.4
$ws.Cells.Item(32, 8).Value = 4833  # H32: 5000 -> 4833
$ws.Cells.Item(32, 9).Value = 4833  # I32: 5000 -> 4833
$ws.Cells.Item(32, 11).Value = 4833  # K32: 5000 -> 4833
$ws.Cells.Item(32, 13).Value = -4516  # M32: -4683 -> -4516
$ws.Cells.Item(33, 8).Value = 26392.166  # H33: 23979 -> 26392.166
$ws.Cells.Item(33, 9).Value = 28005  # I33: 23378.75 -> 28005
$ws.Cells.Item(33, 11).Value = 28005  # K33: 23378.75 -> 28005
$ws.Cells.Item(33, 13).Value = -27715  # M33: -23088.75 -> -27715
$ws.Cells.Item(40, 8).Value = 4441.75  # H40: 4104.6665 -> 4441.75
$ws.Cells.Item(40, 9).Value = 1777  # I40: 3157 -> 1777
$ws.Cells.Item(40, 10).Value = 5330  # J40: 6000 -> 5330
$ws.Cells.Item(40, 11).Value = 1777  # K40: 3157 -> 1777
$ws.Cells.Item(40, 12).Value = 5330  # L40: 6000 -> 5330
$ws.Cells.Item(40, 13).Value = -1641  # M40: -3021 -> -1641
$ws.Cells.Item(40, 14).Value = -5602  # N40: -6272 -> -5602
$ws.Cells.Item(55, 8).Value = 1313.2  # H55: 1313.4 -> 1313.2
$ws.Cells.Item(55, 9).Value = 357.1111  # I55: 357.44446 -> 357.1111
$ws.Cells.Item(55, 11).Value = 357.1111  # K55: 357.44446 -> 357.1111
$ws.Cells.Item(55, 13).Value = -184.1111  # M55: -184.44446 -> -184.1111
$ws.Cells.Item(61, 8).Value = 3233.0557  # H61: 3535.158 -> 3233.0557
$ws.Cells.Item(61, 9).Value = 2705.1333  # I61: 2911.7334 -> 2705.1333
$ws.Cells.Item(61, 10).Value = 5872.6665  # J61: 5873 -> 5872.6665
$ws.Cells.Item(61, 11).Value = 2705.1333  # K61: 2911.7334 -> 2705.1333
$ws.Cells.Item(61, 12).Value = 5872.6665  # L61: 5873 -> 5872.6665
$ws.Cells.Item(61, 13).Value = -2503.1333  # M61: -2709.7334 -> -2503.1333
$ws.Cells.Item(61, 14).Value = -6276.6665  # N61: -6277 -> -6276.6665
$ws.Cells.Item(76, 8).Value = 8166.3335  # H76: 67500 -> 8166.3335
$ws.Cells.Item(76, 10).Value = 8166.3335  # J76: 67500 -> 8166.3335
$ws.Cells.Item(76, 12).Value = 8166.3335  # L76: 67500 -> 8166.3335
$ws.Cells.Item(76, 14).Value = -8842.333500000001  # N76: -68176 -> -8842.333500000001
$ws.Cells.Item(79, 8).Value = 8166.3335  # H79: 67500 -> 8166.3335
$ws.Cells.Item(79, 10).Value = 8166.3335  # J79: 67500 -> 8166.3335
$ws.Cells.Item(79, 12).Value = 8166.3335  # L79: 67500 -> 8166.3335
$ws.Cells.Item(79, 14).Value = -10506.3335  # N79: -69840 -> -10506.3335
$ws.Cells.Item(93, 8).Value = 1514.3103  # H93: 1585.7407 -> 1514.3103
$ws.Cells.Item(93, 9).Value = 1405.5  # I93: 1491.05 -> 1405.5
$ws.Cells.Item(93, 11).Value = 1405.5  # K93: 1491.05 -> 1405.5
$ws.Cells.Item(93, 13).Value = -157.5  # M93: -243.05 -> -157.5
$ws.Cells.Item(113, 8).Value = 3233.0557  # H113: 3535.158 -> 3233.0557
$ws.Cells.Item(113, 9).Value = 2705.1333  # I113: 2911.7334 -> 2705.1333
$ws.Cells.Item(113, 10).Value = 5872.6665  # J113: 5873 -> 5872.6665
$ws.Cells.Item(113, 11).Value = 2705.1333  # K113: 2911.7334 -> 2705.1333
$ws.Cells.Item(113, 12).Value = 5872.6665  # L113: 5873 -> 5872.6665
$ws.Cells.Item(113, 13).Value = -535.1333  # M113: -741.7334000000001 -> -535.1333
$ws.Cells.Item(113, 14).Value = -10212.6665  # N113: -10213 -> -10212.6665
$ws.Cells.Item(122, 8).Value = 4498.4863  # H122: 4540.5 -> 4498.4863
$ws.Cells.Item(122, 9).Value = 3969.1667  # I122: 4027 -> 3969.1667
$ws.Cells.Item(122, 11).Value = 11907.5001  # K122: 12081 -> 11907.5001
$ws.Cells.Item(122, 13).Value = -9457.500100000001  # M122: -9631 -> -9457.500100000001
$ws.Cells.Item(132, 8).Value = 61007.285  # H132: 63867.65 -> 61007.285
$ws.Cells.Item(132, 9).Value = 95343.16  # I132: 102971.75 -> 95343.16
$ws.Cells.Item(132, 11).Value = 286029.48  # K132: 308915.25 -> 286029.48
$ws.Cells.Item(132, 13).Value = -283499.48  # M132: -306385.25 -> -283499.48
$ws.Cells.Item(136, 8).Value = 5500.4287  # H136: 6101.6 -> 5500.4287
$ws.Cells.Item(136, 9).Value = 4282.3335  # I136: 4423.5 -> 4282.3335
$ws.Cells.Item(136, 10).Value = 6414  # J136: 7220.3335 -> 6414
$ws.Cells.Item(136, 11).Value = 12847.0005  # K136: 13270.5 -> 12847.0005
$ws.Cells.Item(136, 12).Value = 19242  # L136: 21661.0005 -> 19242
$ws.Cells.Item(136, 13).Value = -10297.0005  # M136: -10720.5 -> -10297.0005
$ws.Cells.Item(136, 14).Value = -24342  # N136: -26761.0005 -> -24342

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1408.9642  # H113: 1448.2593 -> 1408.9642
$ws.Cells.Item(113, 9).Value = 1252.8  # I113: 1300.421 -> 1252.8
$ws.Cells.Item(113, 11).Value = 3758.4  # K113: 3901.263 -> 3758.4
$ws.Cells.Item(113, 13).Value = -1588.4  # M113: -1731.263 -> -1588.4
$ws.Cells.Item(122, 8).Value = 2104  # H122: 2053.4443 -> 2104
$ws.Cells.Item(122, 9).Value = 1464.6666  # I122: 1510.75 -> 1464.6666
$ws.Cells.Item(122, 11).Value = 4393.9998  # K122: 4532.25 -> 4393.9998
$ws.Cells.Item(122, 13).Value = -1943.9998  # M122: -2082.25 -> -1943.9998
$ws.Cells.Item(136, 8).Value = 2813.449  # H136: 2859.5833 -> 2813.449
$ws.Cells.Item(136, 9).Value = 2375.8838  # I136: 2418.1904 -> 2375.8838
$ws.Cells.Item(136, 11).Value = 7127.651400000001  # K136: 7254.5712 -> 7127.651400000001
$ws.Cells.Item(136, 13).Value = -4577.651400000001  # M136: -4704.5712 -> -4577.651400000001
